# Updated cryptos list on Sun Oct  6 06:33:42 UTC 2024 with GitHub Actions
#
# Refreshes the hourly price/volume snapshot on Sheet1 (crypto ranking
# table). Most rows just get new Price (col D) / Volume(1h) (col E) readings;
# a couple of rows (20/21 and 43/44) swapped ranking order so Coin (col B)
# and Link (col C) are rewritten too.
#
# Numeric-looking Price strings (e.g. "563.88") are written with a leading
# apostrophe so Excel keeps them as literal text (preserving trailing
# zeros/exact formatting) instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.023.65'
$ws.Range("E2").Value = '  -0.22%  '

$ws.Range("D3").Value = '2.420.86'
$ws.Range("E3").Value = '  -0.04%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '''563.88'
$ws.Range("E5").Value = '  +0.92%  '

$ws.Range("D6").Value = '''143.30'
$ws.Range("E6").Value = '  -0.24%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = '''0.532'
$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  +0.21%  '

$ws.Range("E10").Value = '  -0.87%  '

$ws.Range("D11").Value = '''5.22'
$ws.Range("E11").Value = '  -3.52%  '

$ws.Range("D12").Value = '''0.351'
$ws.Range("E12").Value = '  -0.50%  '

$ws.Range("D13").Value = '''26.02'
$ws.Range("E13").Value = '  -0.79%  '

$ws.Range("E14").Value = '  -1.50%  '

$ws.Range("D15").Value = '2.856.88'
$ws.Range("E15").Value = '  +0.29%  '

$ws.Range("D16").Value = '61.883.49'
$ws.Range("E16").Value = '  -0.14%  '

$ws.Range("D17").Value = '2.412.73'
$ws.Range("E17").Value = '  -0.23%  '

$ws.Range("D18").Value = '''11.37'
$ws.Range("E18").Value = '  +1.55%  '

$ws.Range("D19").Value = '''324.46'
$ws.Range("E19").Value = '  +0.10%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '''6.85'
$ws.Range("E20").Value = '  +1.17%  '

$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").Value = '''4.14'
$ws.Range("E21").Value = '  -1.37%  '

$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("D23").Value = '''66.78'
$ws.Range("E23").Value = '  +2.08%  '

$ws.Range("D24").Value = '''1.73'
$ws.Range("E24").Value = '  +0.34%  '

$ws.Range("E25").Value = '  -2.76%  '

$ws.Range("D26").Value = '''556.39'
$ws.Range("E26").Value = '  -6.33%  '

$ws.Range("D27").Value = '2.539.73'
$ws.Range("E27").Value = '  +0.18%  '

$ws.Range("E28").Value = '  +0.36%  '

$ws.Range("D29").Value = '0.0₃0938'
$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("D30").Value = '''8.20'
$ws.Range("E30").Value = '  -1.30%  '

$ws.Range("E31").Value = '  -3.51%  '

$ws.Range("E32").Value = '  -1.87%  '

$ws.Range("D33").Value = '''1.88'
$ws.Range("E33").Value = '  -0.02%  '

$ws.Range("D34").Value = '''1.51'
$ws.Range("E34").Value = '  -3.12%  '

$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("D36").Value = '''4.75'
$ws.Range("E36").Value = '  -0.52%  '

$ws.Range("E37").Value = '  -1.46%  '

$ws.Range("D38").Value = '''154.03'
$ws.Range("E38").Value = '  +1.28%  '

$ws.Range("D39").Value = '''5.44'
$ws.Range("E39").Value = '  -4.69%  '

$ws.Range("D40").Value = '''18.56'
$ws.Range("E40").Value = '  -0.87%  '

$ws.Range("D41").Value = '''1.81'
$ws.Range("E41").Value = '  -0.65%  '

$ws.Range("D42").Value = '''0.991'
$ws.Range("E42").Value = '  -0.84%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '''2.24'
$ws.Range("E43").Value = '  -5.31%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '''147.14'
$ws.Range("E44").Value = '  -2.52%  '

$ws.Range("E45").Value = '  -0.31%  '

$ws.Range("D46").Value = '''0.0528'
$ws.Range("E46").Value = '  -2.03%  '

$ws.Range("D47").Value = '''19.86'
$ws.Range("E47").Value = '  -2.02%  '

$ws.Range("D48").Value = '''0.593'
$ws.Range("E48").Value = '  +0.10%  '

$ws.Range("D49").Value = '''0.0921'
$ws.Range("E49").Value = '  +0.17%  '

$ws.Range("D50").Value = '''0.0227'
$ws.Range("E50").Value = '  -0.88%  '

$ws.Range("E51").Value = '  +0.65%  '
